$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, "28.676.77", "  -0.59%  "),
    @(3, "1.884.59", "  +0.18%  "),
    @(4, "1.007", "  +0.07%  "),
    @(5, "326.46", "  +0.07%  "),
    @(6, "1.006", "  +0.02%  "),
    @(7, "0.4600", "  -1.62%  "),
    @(8, "0.3867", "  -2.05%  "),
    @(9, "46.72", "  +0.22%  "),
    @(10, "0.07881", "  -0.69%  "),
    @(11, "1.003", "  +2.18%  "),
    @(12, "21.77", "  -2.59%  "),
    @(13, "1.876.21", "  -1.04%  "),
    @(14, "7.078", "  +0.78%  "),
    @(15, "5.713", "  -0.68%  "),
    @(16, "0.06988", "  +0.20%  "),
    @(17, "87.61", "  -1.22%  "),
    @(18, "1.007", "  +0.04%  "),
    @(19, "0.00001005", "  -0.47%  "),
    @(20, "17.24", "  +1.34%  "),
    @(21, "1.007", "  +0.23%  "),
    @(22, "28.667.75", "  -0.67%  "),
    @(23, "5.334", "  -0.64%  "),
    @(24, "10.99", "  -1.08%  "),
    @(25, "2.108.96", "  -0.49%  "),
    @(26, "2.054", "  -3.27%  "),
    @(27, "154.54", "  +0.56%  "),
    @(28, "19.56", "  +0.88%  "),
    @(29, "5.856", "  +1.47%  "),
    @(30, "1.967", "  -2.02%  "),
    @(31, "118.65", "  -1.18%  "),
    @(32, "0.09348", "  -0.58%  "),
    @(33, "0.9239", "  -1.95%  "),
    @(34, "5.327", "  +0.06%  "),
    @(35, "1.339", "  -1.33%  "),
    @(36, "3.267", "  -2.61%  "),
    @(37, "0.05770", "  -2.57%  "),
    @(38, "1.166", "  +1.26%  "),
    @(39, "8.017", "  +1.20%  "),
    @(40, "0.02074", "  -2.32%  "),
    @(41, "0.5689", "  -0.72%  "),
    @(42, "0.1797", "  -0.11%  "),
    @(43, "9.762", "  -2.48%  "),
    @(44, "11.85", "  -0.20%  "),
    @(45, "0.5377", "  +0.53%  "),
    @(46, "0.07151", "  -1.74%  "),
    @(47, "2.152", "  +1.33%  "),
    @(48, "1.846", "  -0.36%  "),
    @(49, "1.116", "  -3.00%  "),
    @(50, "112.49", "  -1.79%  "),
    @(51, "2.485", "  +4.84%  "),
)

foreach ($item in $updates) {
    $row = $item[0]
    $dVal = $item[1]
    $eVal = $item[2]
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = "'" + $dVal
    $dCell.Style = "Normal"
    $ws.Cells.Item($row, 5).Value = $eVal
}
